$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 681, shifting existing rows 681:718 down to 682:719
$ws.Rows.Item(681).Insert()

# Populate the newly inserted row 681 with the new weekly data point
$ws.Range("A681").Value = 8
$ws.Range("B681").Value = "Terminal La Palmera de La Serena"
$ws.Range("C681").Value = "Coquimbo"
$ws.Range("D681").Value = 45267
$ws.Range("E681").Value = 4
$ws.Range("F681").Value = 100112043
$ws.Range("G681").Value = "Pepino dulce"
$ws.Range("H681").Value = "Sin especificar"
$ws.Range("I681").Value = "Segunda"
$ws.Range("J681").Value = 240
$ws.Range("K681").Value = 19000
$ws.Range("L681").Value = 20000
$ws.Range("M681").Value = 19500
$ws.Range("N681").Value = "`$/bandeja 18 kilos"
$ws.Range("O681").Value = "Provincia de Limarí"
$ws.Range("P681").Value = 1083
$ws.Range("Q681").Value = 18
$ws.Range("R681").Value = "Hortaliza"
